$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("forecasts (all origional)")

# ---------------------------------------------------------------------------
# 1) Re-run the forecast error metric: RMSE -> MAPE on the
#    "forecasts (all origional)" sheet (AB:AL helper columns + AN summary).
# ---------------------------------------------------------------------------

# Row 2 holds explicit (non-shared) formulas.
$ws3.Range("AB2").Formula = "=ABS((P2-O2)/O2*100)"
$ws3.Range("AC2").Formula = "=ABS((Q3-O3)/O3*100)"
$ws3.Range("AD2").Formula = "=ABS((R4-O4)/O4*100)"
$ws3.Range("AE2").Formula = "=ABS((S5-O5)/O5*100)"
$ws3.Range("AF2").Formula = "=ABS((T6-O6)/O6*100)"
$ws3.Range("AG2").Formula = "=ABS((U7-O7)/O7*100)"
$ws3.Range("AH2").Formula = "=ABS((V8-O8)/O8*100)"
$ws3.Range("AI2").Formula = "=ABS((W9-O9)/O9*100)"
$ws3.Range("AJ2").Formula = "=ABS((X10-O10)/O10*100)"
$ws3.Range("AK2").Formula = "=ABS((Y11-O11)/O11*100)"
$ws3.Range("AL2").Formula = "=ABS((Z12-O12)/O12*100)"

# Rows 3:13 - set the whole column range at once so the engine fills the
# relative references down and keeps them as one shared-formula group.
$ws3.Range("AB3:AB13").Formula = "=ABS((P3-O3)/O3*100)"
$ws3.Range("AC3:AC13").Formula = "=ABS((Q4-O4)/O4*100)"
$ws3.Range("AD3:AD13").Formula = "=ABS((R5-O5)/O5*100)"
$ws3.Range("AE3:AE13").Formula = "=ABS((S6-O6)/O6*100)"
$ws3.Range("AF3:AF13").Formula = "=ABS((T7-O7)/O7*100)"
$ws3.Range("AG3:AG13").Formula = "=ABS((U8-O8)/O8*100)"
$ws3.Range("AH3:AH13").Formula = "=ABS((V9-O9)/O9*100)"
$ws3.Range("AI3:AI13").Formula = "=ABS((W10-O10)/O10*100)"
$ws3.Range("AJ3:AJ13").Formula = "=ABS((X11-O11)/O11*100)"
$ws3.Range("AK3:AK13").Formula = "=ABS((Y12-O12)/O12*100)"
$ws3.Range("AL3:AL13").Formula = "=ABS((Z13-O13)/O13*100)"

# AN column: the per-h summary switches from SQRT(AVERAGE(...)) (rmse) to
# plain AVERAGE(...) (mape).
$ws3.Range("AN2").Formula = "=AVERAGE(AB2:AL2)"
$ws3.Range("AN3:AN13").Formula = "=AVERAGE(AB3:AL3)"

# ---------------------------------------------------------------------------
# 2) Drop the stale _xlchart.* defined names (leftover chart bookmarks that
#    no longer correspond to anything in the workbook).
# ---------------------------------------------------------------------------
for ($i = $wb.Names.Count; $i -ge 1; $i--) {
    $wb.Names.Item($i).Delete()
}

# ---------------------------------------------------------------------------
# 3) Add the new "mape" presentation sheet after "rmse", mirroring the
#    existing "rmse" sheet layout (h / mape columns) but sourced from the
#    freshly recomputed AN column above.
# ---------------------------------------------------------------------------
$rmseSheet = $wb.Worksheets.Item("rmse")
$mapeSheet = $wb.Worksheets.Add($null, $rmseSheet)
$mapeSheet.Name = "mape"

$mapeSheet.Range("A1").Value = "h"
$mapeSheet.Range("B1").Value = "mape"
$mapeSheet.Range("B1").Font.Bold = $true

for ($h = 1; $h -le 12; $h++) {
    $row = $h + 1
    $mapeSheet.Cells.Item($row, 1).Value = $h
    $mapeSheet.Cells.Item($row, 2).Value = $ws3.Cells.Item($h + 1, 40).Value2
}

# Make "mape" the active sheet/tab and match the recorded selection.
$mapeSheet.Activate()
$mapeSheet.Range("F11").Select()
